# Add a new "UK" test-data sheet, modeled on the existing "Portugal" sheet
# (same column widths, same 21-row "P panel" layout), but place it using
# "Poland" as the structural base so that row heights / wrapping stay at
# their defaults (matching the target workbook), then insert the two rows
# that Poland is missing (P32AR / P32DR) to reach the full 21-row layout.

$wb = $excel.ActiveWorkbook

# Poland is currently the last sheet in the workbook; duplicate it to the
# end of the workbook to serve as the new "UK" sheet.
$template = $wb.Worksheets.Item("Poland")
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "UK"

# Poland's repeater list is missing "P32AR" / "P32DR" rows that the other
# 21-row market sheets (e.g. Portugal) have. Insert two rows above the
# current row 16 ("PR1DS") to make room for them.
$newSheet.Rows("16:17").Insert()

# The inserted rows come back with no explicit style; copy the (already
# correct) formatting from row 18 (pushed-down "PR1DS" row, style index 3)
# onto the two new rows so no new style entries are created.
$newSheet.Range("A18").Copy()
$newSheet.Range("A16:A17").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$newSheet.Range("A16").Value = "P32AR"
$newSheet.Range("A17").Value = "P32DR"

# B4 ("User Story" value cell) in Poland has no explicit style, but the
# target sheet carries the bordered style (index 3) used elsewhere in the
# column. Copy that formatting over before setting the value.
$newSheet.Range("A16").Copy()
$newSheet.Range("B4").PasteSpecial(-4122)        # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the market-specific values.
$newSheet.Range("B4").Value = "NGC-2741/T3364"
$newSheet.Range("B2").Value = "UK Market"

# Leave B4 selected/active, matching the saved selection state of the
# authored sheet.
$newSheet.Range("B4").Select()
